$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 74 values (quarter 01-01-2021)
$ws.Range("B74").Value = 2125
$ws.Range("E74").Value = 2041
$ws.Range("F74").Value = 742
$ws.Range("G74").Value = 1299

# Add new row 75 (quarter 01-04-2021)
# Use a formula-then-paste-values trick so the date-like text "01-04-2021"
# is stored as a plain shared string (not auto-converted to a date serial).
$ws.Range("A75").Formula = "=""01-04-2021"""
$ws.Range("A75").Copy()
$ws.Range("A75").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("B75").Value = 2081
$ws.Range("C75").Value = 107
$ws.Range("D75").Value = 107
$ws.Range("E75").Value = 1974
$ws.Range("F75").Value = 683
$ws.Range("G75").Value = 1291
